$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.184.56'
$ws.Range("E2").Value = '  +1.66%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.800.50'
$ws.Range("E3").Value = '  +2.39%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  -0.32%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '324.52'
$ws.Range("E5").Value = '  -0.68%  '
$ws.Range("E6").Value = '  -0.11%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4293'
$ws.Range("E7").Value = '  -3.09%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3636'
$ws.Range("E8").Value = '  -2.83%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '44.70'
$ws.Range("E9").Value = '  -1.68%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07563'
$ws.Range("E10").Value = '  -1.39%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.122'
$ws.Range("E11").Value = '  -0.19%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.001'
$ws.Range("E12").Value = '  -0.14%  '
$ws.Range("E13").Value = '  +0.26%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.187'
$ws.Range("E14").Value = '  -0.15%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.380'
$ws.Range("E15").Value = '  -0.86%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.814.31'
$ws.Range("E16").Value = '  +3.24%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '92.94'
$ws.Range("E17").Value = '  +4.46%  '
$ws.Range("E18").Value = '  -0.52%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06351'
$ws.Range("E19").Value = '  +1.97%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.000'
$ws.Range("E20").Value = '  -0.06%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.27'
$ws.Range("E21").Value = '  -0.31%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.003'
$ws.Range("E22").Value = '  -2.86%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '28.207.80'
$ws.Range("E23").Value = '  +1.63%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.43'
$ws.Range("E24").Value = '  -1.79%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.178'
$ws.Range("E25").Value = '  -5.90%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '159.41'
$ws.Range("E26").Value = '  +3.85%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.44'
$ws.Range("E27").Value = '  -1.62%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.018.03'
$ws.Range("E28").Value = '  +3.16%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.243'
$ws.Range("E29").Value = '  -4.86%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '128.11'
$ws.Range("E30").Value = '  -0.23%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.181'
$ws.Range("E31").Value = '  -2.80%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.888'
$ws.Range("E32").Value = '  +2.19%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.09035'
$ws.Range("E33").Value = '  -3.19%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.533'
$ws.Range("E34").Value = '  -3.17%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '12.81'
$ws.Range("E35").Value = '  +0.68%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02362'
$ws.Range("E36").Value = '  +1.78%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.141'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.6522'
$ws.Range("E38").Value = '  +0.42%  '
$ws.Range("B39").Value = 'Hedera'
$ws.Range("C39").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06128'
$ws.Range("E39").Value = '  -0.26%  '
$ws.Range("B40").Value = 'Algorand'
$ws.Range("C40").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.2126'
$ws.Range("E40").Value = '  -2.76%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.199'
$ws.Range("E41").Value = '  -0.22%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.431'
$ws.Range("E42").Value = '  +0.91%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '7.961'
$ws.Range("E43").Value = '  -0.37%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.9992'
$ws.Range("E44").Value = '  -0.15%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.67'
$ws.Range("E45").Value = '  -0.75%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.6035'
$ws.Range("E46").Value = '  +0.20%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.713'
$ws.Range("E47").Value = '  -1.41%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '125.64'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.999'
$ws.Range("E49").Value = '  +0.01%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.159'
$ws.Range("E50").Value = '  +1.89%  '
$ws.Range("E51").Value = '  +1.05%  '
